$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions scheduled update).
# Several "Price" column values are plain decimal numbers (e.g. "0.999") which
# Excel would otherwise auto-convert to the Number type on assignment, but the
# source data must stay text (consistent with values like "41.586.41" that can
# never be numeric), so those specific cells are pre-formatted as Text first.
$textCells = $ws.Range("D4,D5,D6,D9,D10,D14,D15,D17,D20,D21,D22,D23,D26,D27,D30,D31,D32,D35,D36,D37,D39,D41,D45,D46,D49,D50,D51")
$textCells.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "41.586.41"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.467.72"
$ws.Range("E3").Value = "  -0.60%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.19%  "

# Row 5 - BNB
$ws.Range("D5").Value = "317.58"
$ws.Range("E5").Value = "  +1.31%  "

# Row 6 - Solana
$ws.Range("D6").Value = "92.03"
$ws.Range("E6").Value = "  -0.42%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.64%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.13%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.515"

# Row 10 - Avalanche
$ws.Range("D10").Value = "32.87"
$ws.Range("E10").Value = "  -0.08%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +7.42%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.37%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.847.05"
$ws.Range("E13").Value = "  -0.59%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "6.87"
$ws.Range("E14").Value = "  -0.72%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "15.53"
$ws.Range("E15").Value = "  -5.22%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.458.29"
$ws.Range("E16").Value = "  -0.34%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.790"
$ws.Range("E17").Value = "  +1.79%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "41.537.32"
$ws.Range("E18").Value = "  -0.04%  "

# Row 19 - Uniswap->ShibaInu (swap)
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0948"
$ws.Range("E19").Value = "  +0.14%  "

# Row 20 - ShibaInu->Uniswap (swap)
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "6.44"
$ws.Range("E20").Value = "  -1.95%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "71.09"
$ws.Range("E21").Value = "  -2.13%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("D22").Value = "11.30"
$ws.Range("E22").Value = "  +0.68%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "239.10"
$ws.Range("E23").Value = "  +0.77%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  +0.36%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  +0.97%  "

# Row 26 - Dai
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.19%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "24.56"

# Row 28 - Toncoin
$ws.Range("E28").Value = "  +2.55%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +1.06%  "

# Row 30 - InjectiveProtocol
$ws.Range("D30").Value = "36.04"
$ws.Range("E30").Value = "  +0.32%  "

# Row 31 - Monero
$ws.Range("D31").Value = "161.22"
$ws.Range("E31").Value = "  +2.05%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "5.50"
$ws.Range("E32").Value = "  +0.44%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  +0.07%  "

# Row 34 - WEMIXToken
$ws.Range("E34").Value = "  +0.35%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "0.0765"
$ws.Range("E35").Value = "  +0.59%  "

# Row 36 - Celestia
$ws.Range("D36").Value = "17.24"
$ws.Range("E36").Value = "  -0.81%  "

# Row 37 - LidoDAOToken
$ws.Range("D37").Value = "2.90"
$ws.Range("E37").Value = "  -0.94%  "

# Row 38 - Stellar
$ws.Range("E38").Value = "  +1.58%  "

# Row 39 - ARBITRUM
$ws.Range("D39").Value = "1.83"
$ws.Range("E39").Value = "  -0.51%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  -3.39%  "

# Row 41 - RenderToken
$ws.Range("D41").Value = "3.97"
$ws.Range("E41").Value = "  -2.83%  "

# Row 42 - ApeXProtocol
$ws.Range("E42").Value = "  +3.25%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.983.30"
$ws.Range("E43").Value = "  +0.50%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -0.15%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "18.79"
$ws.Range("E45").Value = "  -2.25%  "

# Row 46 - NEARProtocol
$ws.Range("D46").Value = "2.98"
$ws.Range("E46").Value = "  +0.57%  "

# Row 47 - FraxShare
$ws.Range("E47").Value = "  +2.57%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "2.705.80"
$ws.Range("E48").Value = "  -0.55%  "

# Row 49 - Aave
$ws.Range("D49").Value = "97.26"
$ws.Range("E49").Value = "  -0.91%  "

# Row 50 - BitcoinSV
$ws.Range("D50").Value = "74.09"
$ws.Range("E50").Value = "  +2.02%  "

# Row 51 - ordi
$ws.Range("D51").Value = "67.14"
$ws.Range("E51").Value = "  -2.02%  "
